$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = [double]"32.867558"
$ws.Range("H2").Value2 = [double]"98.60267399999999"
$ws.Range("I2").Value2 = [double]"0.5341839962912384"
$ws.Range("J2").Value2 = [double]"0.5341839962912384"
$ws.Range("M2").Value2 = [double]"32.867558"
$ws.Range("N2").Value2 = [double]"98.60267399999999"
$ws.Range("O2").Value2 = [double]"0.5341839962912384"
$ws.Range("P2").Value2 = [double]"0.5341839962912384"
$ws.Range("Q2").Value2 = [double]"1080.276368883364"
$ws.Range("R2").Value2 = [double]"9722.487319950274"
$ws.Range("S2").Value2 = [double]"0.2853525418936778"
$ws.Range("T2").Value2 = [double]"0.2853525418936778"
$ws.Range("G3").Value2 = [double]"32.867558"
$ws.Range("H3").Value2 = [double]"98.60267399999999"
$ws.Range("I3").Value2 = [double]"0.5341839962912384"
$ws.Range("J3").Value2 = [double]"0.5341839962912384"
$ws.Range("O3").Value2 = [double]"0.1136291069637658"
$ws.Range("P3").Value2 = [double]"0.1136291069637658"
$ws.Range("Q3").Value2 = [double]"229.7913077189086"
$ws.Range("R3").Value2 = [double]"2068.121769470178"
$ws.Range("S3").Value2 = [double]"0.06069885045290901"
$ws.Range("T3").Value2 = [double]"0.06069885045290901"
$ws.Range("G4").Value2 = [double]"32.867558"
$ws.Range("H4").Value2 = [double]"98.60267399999999"
$ws.Range("I4").Value2 = [double]"0.5341839962912384"
$ws.Range("J4").Value2 = [double]"0.5341839962912384"
$ws.Range("M4").Value2 = [double]"21.63740966666667"
$ws.Range("N4").Value2 = [double]"64.912229"
$ws.Range("O4").Value2 = [double]"0.3516646404071356"
$ws.Range("P4").Value2 = [double]"0.3516646404071356"
$ws.Range("Q4").Value2 = [double]"711.1688171889273"
$ws.Range("R4").Value2 = [double]"6400.519354700346"
$ws.Range("S4").Value2 = [double]"0.187853622967005"
$ws.Range("T4").Value2 = [double]"0.187853622967005"
$ws.Range("G5").Value2 = [double]"32.867558"
$ws.Range("H5").Value2 = [double]"98.60267399999999"
$ws.Range("I5").Value2 = [double]"0.5341839962912384"
$ws.Range("J5").Value2 = [double]"0.5341839962912384"
$ws.Range("M5").Value2 = [double]"0.03213366666666666"
$ws.Range("N5").Value2 = [double]"0.096401"
$ws.Range("O5").Value2 = [double]"0.0005222563378602863"
$ws.Range("P5").Value2 = [double]"0.0005222563378602864"
$ws.Range("Q5").Value2 = [double]"1.056155152919333"
$ws.Range("R5").Value2 = [double]"9.505396376274"
$ws.Range("S5").Value2 = [double]"0.0002789809776466349"
$ws.Range("T5").Value2 = [double]"0.000278980977646635"
$ws.Range("I6").Value2 = [double]"0.1136291069637658"
$ws.Range("J6").Value2 = [double]"0.1136291069637658"
$ws.Range("M6").Value2 = [double]"32.867558"
$ws.Range("N6").Value2 = [double]"98.60267399999999"
$ws.Range("O6").Value2 = [double]"0.5341839962912384"
$ws.Range("P6").Value2 = [double]"0.5341839962912384"
$ws.Range("Q6").Value2 = [double]"229.7913077189086"
$ws.Range("R6").Value2 = [double]"2068.121769470178"
$ws.Range("S6").Value2 = [double]"0.06069885045290901"
$ws.Range("T6").Value2 = [double]"0.06069885045290901"
$ws.Range("I7").Value2 = [double]"0.1136291069637658"
$ws.Range("J7").Value2 = [double]"0.1136291069637658"
$ws.Range("O7").Value2 = [double]"0.1136291069637658"
$ws.Range("P7").Value2 = [double]"0.1136291069637658"
$ws.Range("S7").Value2 = [double]"0.01291157394938293"
$ws.Range("T7").Value2 = [double]"0.01291157394938294"
$ws.Range("I8").Value2 = [double]"0.1136291069637658"
$ws.Range("J8").Value2 = [double]"0.1136291069637658"
$ws.Range("M8").Value2 = [double]"21.63740966666667"
$ws.Range("N8").Value2 = [double]"64.912229"
$ws.Range("O8").Value2 = [double]"0.3516646404071356"
$ws.Range("P8").Value2 = [double]"0.3516646404071356"
$ws.Range("Q8").Value2 = [double]"151.2764855531125"
$ws.Range("R8").Value2 = [double]"1361.488369978013"
$ws.Range("S8").Value2 = [double]"0.03995933904019666"
$ws.Range("T8").Value2 = [double]"0.03995933904019666"
$ws.Range("I9").Value2 = [double]"0.1136291069637658"
$ws.Range("J9").Value2 = [double]"0.1136291069637658"
$ws.Range("M9").Value2 = [double]"0.03213366666666666"
$ws.Range("N9").Value2 = [double]"0.096401"
$ws.Range("O9").Value2 = [double]"0.0005222563378602863"
$ws.Range("P9").Value2 = [double]"0.0005222563378602864"
$ws.Range("Q9").Value2 = [double]"0.2246603561218889"
$ws.Range("R9").Value2 = [double]"2.021943205097"
$ws.Range("S9").Value2 = [double]"5.934352127723109E-05"
$ws.Range("T9").Value2 = [double]"5.934352127723111E-05"
$ws.Range("G10").Value2 = [double]"21.63740966666667"
$ws.Range("H10").Value2 = [double]"64.912229"
$ws.Range("I10").Value2 = [double]"0.3516646404071356"
$ws.Range("J10").Value2 = [double]"0.3516646404071356"
$ws.Range("M10").Value2 = [double]"32.867558"
$ws.Range("N10").Value2 = [double]"98.60267399999999"
$ws.Range("O10").Value2 = [double]"0.5341839962912384"
$ws.Range("P10").Value2 = [double]"0.5341839962912384"
$ws.Range("Q10").Value2 = [double]"711.1688171889273"
$ws.Range("R10").Value2 = [double]"6400.519354700346"
$ws.Range("S10").Value2 = [double]"0.187853622967005"
$ws.Range("T10").Value2 = [double]"0.187853622967005"
$ws.Range("G11").Value2 = [double]"21.63740966666667"
$ws.Range("H11").Value2 = [double]"64.912229"
$ws.Range("I11").Value2 = [double]"0.3516646404071356"
$ws.Range("J11").Value2 = [double]"0.3516646404071356"
$ws.Range("O11").Value2 = [double]"0.1136291069637658"
$ws.Range("P11").Value2 = [double]"0.1136291069637658"
$ws.Range("Q11").Value2 = [double]"151.2764855531125"
$ws.Range("R11").Value2 = [double]"1361.488369978013"
$ws.Range("S11").Value2 = [double]"0.03995933904019666"
$ws.Range("T11").Value2 = [double]"0.03995933904019666"
$ws.Range("G12").Value2 = [double]"21.63740966666667"
$ws.Range("H12").Value2 = [double]"64.912229"
$ws.Range("I12").Value2 = [double]"0.3516646404071356"
$ws.Range("J12").Value2 = [double]"0.3516646404071356"
$ws.Range("M12").Value2 = [double]"21.63740966666667"
$ws.Range("N12").Value2 = [double]"64.912229"
$ws.Range("O12").Value2 = [double]"0.3516646404071356"
$ws.Range("P12").Value2 = [double]"0.3516646404071356"
$ws.Range("Q12").Value2 = [double]"468.1774970831601"
$ws.Range("R12").Value2 = [double]"4213.59747374844"
$ws.Range("S12").Value2 = [double]"0.12366801931268"
$ws.Range("T12").Value2 = [double]"0.12366801931268"
$ws.Range("G13").Value2 = [double]"21.63740966666667"
$ws.Range("H13").Value2 = [double]"64.912229"
$ws.Range("I13").Value2 = [double]"0.3516646404071356"
$ws.Range("J13").Value2 = [double]"0.3516646404071356"
$ws.Range("M13").Value2 = [double]"0.03213366666666666"
$ws.Range("N13").Value2 = [double]"0.096401"
$ws.Range("O13").Value2 = [double]"0.0005222563378602863"
$ws.Range("P13").Value2 = [double]"0.0005222563378602864"
$ws.Range("Q13").Value2 = [double]"0.6952893097587778"
$ws.Range("R13").Value2 = [double]"6.257603787829"
$ws.Range("S13").Value2 = [double]"0.0001836590872539851"
$ws.Range("T13").Value2 = [double]"0.0001836590872539851"
$ws.Range("G14").Value2 = [double]"0.03213366666666666"
$ws.Range("H14").Value2 = [double]"0.096401"
$ws.Range("I14").Value2 = [double]"0.0005222563378602863"
$ws.Range("J14").Value2 = [double]"0.0005222563378602864"
$ws.Range("M14").Value2 = [double]"32.867558"
$ws.Range("N14").Value2 = [double]"98.60267399999999"
$ws.Range("O14").Value2 = [double]"0.5341839962912384"
$ws.Range("P14").Value2 = [double]"0.5341839962912384"
$ws.Range("Q14").Value2 = [double]"1.056155152919333"
$ws.Range("R14").Value2 = [double]"9.505396376274"
$ws.Range("S14").Value2 = [double]"0.0002789809776466349"
$ws.Range("T14").Value2 = [double]"0.000278980977646635"
$ws.Range("G15").Value2 = [double]"0.03213366666666666"
$ws.Range("H15").Value2 = [double]"0.096401"
$ws.Range("I15").Value2 = [double]"0.0005222563378602863"
$ws.Range("J15").Value2 = [double]"0.0005222563378602864"
$ws.Range("O15").Value2 = [double]"0.1136291069637658"
$ws.Range("P15").Value2 = [double]"0.1136291069637658"
$ws.Range("Q15").Value2 = [double]"0.2246603561218889"
$ws.Range("R15").Value2 = [double]"2.021943205097"
$ws.Range("S15").Value2 = [double]"5.934352127723109E-05"
$ws.Range("T15").Value2 = [double]"5.934352127723111E-05"
$ws.Range("G16").Value2 = [double]"0.03213366666666666"
$ws.Range("H16").Value2 = [double]"0.096401"
$ws.Range("I16").Value2 = [double]"0.0005222563378602863"
$ws.Range("J16").Value2 = [double]"0.0005222563378602864"
$ws.Range("M16").Value2 = [double]"21.63740966666667"
$ws.Range("N16").Value2 = [double]"64.912229"
$ws.Range("O16").Value2 = [double]"0.3516646404071356"
$ws.Range("P16").Value2 = [double]"0.3516646404071356"
$ws.Range("Q16").Value2 = [double]"0.6952893097587778"
$ws.Range("R16").Value2 = [double]"6.257603787829"
$ws.Range("S16").Value2 = [double]"0.0001836590872539851"
$ws.Range("T16").Value2 = [double]"0.0001836590872539851"
$ws.Range("G17").Value2 = [double]"0.03213366666666666"
$ws.Range("H17").Value2 = [double]"0.096401"
$ws.Range("I17").Value2 = [double]"0.0005222563378602863"
$ws.Range("J17").Value2 = [double]"0.0005222563378602864"
$ws.Range("M17").Value2 = [double]"0.03213366666666666"
$ws.Range("N17").Value2 = [double]"0.096401"
$ws.Range("O17").Value2 = [double]"0.0005222563378602863"
$ws.Range("P17").Value2 = [double]"0.0005222563378602864"
$ws.Range("Q17").Value2 = [double]"0.001032572533444444"
$ws.Range("R17").Value2 = [double]"0.009293152801"
$ws.Range("S17").Value2 = [double]"2.727516824352375E-07"
$ws.Range("T17").Value2 = [double]"2.727516824352376E-07"
